# "Generate Report for Handback"
#
# Populates the "Latest Target File" (F) and "Latest Handback File" (G)
# columns on the per-language report sheets (zh-cn, de-de) with hyperlinked
# file names, refreshes the "Status" text workbook-wide (Overview + both
# language sheets) from "Ready for handoff" to the handed-back state, and
# stamps the per-language "Latest Handback DateTime" column with the actual
# handback timestamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# 1. Status: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (Overview!B/C2:3 and the Status column on both language sheets all
#    shared the old string, so every occurrence is refreshed together.)
# ---------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. Latest Handback DateTime (column H) -- real timestamps instead of
#    the zero-date placeholder, per language.
# ---------------------------------------------------------------------
$wsZhCn.Range("H2").Value = "2016-03-22 14:34:33"
$wsZhCn.Range("H3").Value = "2016-03-22 14:34:33"

$wsDeDe.Range("H2").Value = "2016-03-22 14:34:42"
$wsDeDe.Range("H3").Value = "2016-03-22 14:34:42"

# ---------------------------------------------------------------------
# 3. Latest Target File (F) / Latest Handback File (G): new hyperlinked
#    cells on rows 2 and 3 of each language sheet. Existing hyperlinks
#    (A2, D2, A3, D3) are re-created too so the final hyperlink order
#    matches a fresh left-to-right, top-to-bottom regeneration of the
#    sheet (A2, D2, F2, G2, A3, D3, F3, G3).
# ---------------------------------------------------------------------
function Set-ReportSheetHandback {
    param($ws, $xlfName)

    # Capture the existing hyperlink targets (keyed by cell address) before
    # wiping the collection, so A2/D2/A3/D3 can be restored unchanged.
    $existing = @{}
    foreach ($h in $ws.Hyperlinks) {
        $existing[$h.Range.Address()] = $h.Address
    }

    $aTarget = $existing["`$A`$2"]
    $bTarget = $existing["`$A`$3"]
    $dTarget = $existing["`$D`$2"]

    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $aTarget, "", "", "a.md")
    $ws.Hyperlinks.Add($ws.Range("D2"), $dTarget, "", "", $xlfName)
    $ws.Hyperlinks.Add($ws.Range("F2"), $aTarget, "", "", "a.md")
    $ws.Hyperlinks.Add($ws.Range("G2"), $dTarget, "", "", $xlfName)

    $ws.Hyperlinks.Add($ws.Range("A3"), $bTarget, "", "", "b.md")
    $ws.Hyperlinks.Add($ws.Range("D3"), $dTarget, "", "", $xlfName)
    $ws.Hyperlinks.Add($ws.Range("F3"), $aTarget, "", "", "a.md")
    $ws.Hyperlinks.Add($ws.Range("G3"), $dTarget, "", "", $xlfName)

    $ws.Range("F2").Style = "HyperLink"
    $ws.Range("G2").Style = "HyperLink"
    $ws.Range("F3").Style = "HyperLink"
    $ws.Range("G3").Style = "HyperLink"
}

Set-ReportSheetHandback $wsZhCn "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
Set-ReportSheetHandback $wsDeDe "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
